# DOMA-746 add executor & assignee excel export mappers
#
# Rename the `ticket` list accessor to `tickets` in every export-template
# placeholder (rows 2 & 3 of the sheet hold the `{d.ticket[...]....}` /
# `{d.ticket[i + 1]....}` mustache-style tokens), and move the active
# selection to F19 to match the author's saved cursor position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2 and 3, columns A-G (1-7) hold the `{d.ticket[...]. ...}` /
# `{d.ticket[i + 1]. ...}` placeholder strings - rename the `ticket`
# collection to `tickets` wherever it appears.
for ($row = 2; $row -le 3; $row++) {
    for ($col = 1; $col -le 7; $col++) {
        $cell = $ws.Cells.Item($row, $col)
        $current = $cell.Value()
        if ($current -and $current -match "d\.ticket\[") {
            $cell.Value = $current -replace "d\.ticket\[", "d.tickets["
        }
    }
}

# Move/save the active selection at F19 (matches the saved cursor position
# recorded in the sheet view).
$ws.Range("F19").Select()
